# Scheduled-runner update: refresh market-board price snapshots and
# recomputed leve profit figures (currentAveragePrice* / LevePrice* /
# LeveProfit* columns) across the per-job Atomos_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 3120.5134
$ws.Range("I116").Value = 2441.7778
$ws.Range("J116").Value = 4953.1
$ws.Range("K116").Value = 2441.7778
$ws.Range("L116").Value = 4953.1
$ws.Range("M116").Value = 1000.2222
$ws.Range("N116").Value = -11837.1

# Row 130
$ws.Range("H130").Value = 10000
$ws.Range("J130").Value = 10000
$ws.Range("L130").Value = 10000
$ws.Range("N130").Value = -20040

# Row 134
$ws.Range("H134").Value = 29500
$ws.Range("J134").Value = 29500
$ws.Range("L134").Value = 29500
$ws.Range("N134").Value = -39640

# Row 137
$ws.Range("H137").Value = 4549945.5
$ws.Range("I137").Value = 6671860
$ws.Range("K137").Value = 20015580
$ws.Range("M137").Value = -20013030

# Row 139
$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280

# Row 140
$ws.Range("H140").Value = 24572.857
$ws.Range("J140").Value = 24572.857
$ws.Range("L140").Value = 24572.857
$ws.Range("N140").Value = -34932.857

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6480.3228
$ws.Range("I32").Value = 3903.7754
$ws.Range("K32").Value = 3903.7754
$ws.Range("M32").Value = -3616.7754

# Row 35
$ws.Range("H35").Value = 2068.5
$ws.Range("I35").Value = 2068.5
$ws.Range("K35").Value = 2068.5
$ws.Range("M35").Value = -1662.5

# Row 61
$ws.Range("H61").Value = 1825.2941
$ws.Range("I61").Value = 1314.375
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 1314.375
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -1102.375
$ws.Range("N61").Value = -10424

# Row 74
$ws.Range("H74").Value = 637.6667
$ws.Range("I74").Value = 536.125
$ws.Range("J74").Value = 1450
$ws.Range("K74").Value = 536.125
$ws.Range("L74").Value = 1450
$ws.Range("M74").Value = 337.875
$ws.Range("N74").Value = -3198

# Row 77
$ws.Range("H77").Value = 637.6667
$ws.Range("I77").Value = 536.125
$ws.Range("J77").Value = 1450
$ws.Range("K77").Value = 2680.625
$ws.Range("L77").Value = 7250
$ws.Range("M77").Value = 1687.375
$ws.Range("N77").Value = -15986

# Row 122
$ws.Range("H122").Value = 2681.0454
$ws.Range("I122").Value = 1732
$ws.Range("J122").Value = 3819.9
$ws.Range("K122").Value = 5196
$ws.Range("L122").Value = 11459.7
$ws.Range("M122").Value = -2746
$ws.Range("N122").Value = -16359.7

# Row 136
$ws.Range("H136").Value = 1825.2941
$ws.Range("I136").Value = 1314.375
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 3943.125
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -1393.125
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
# Row 58
$ws.Range("H58").Value = 28560
$ws.Range("J58").Value = 28560
$ws.Range("L58").Value = 28560
$ws.Range("N58").Value = -29148

# Row 86
$ws.Range("H86").Value = 1934.9412
$ws.Range("I86").Value = 1283.3334
$ws.Range("K86").Value = 1283.3334
$ws.Range("M86").Value = -160.3334

# Row 89
$ws.Range("H89").Value = 1934.9412
$ws.Range("I89").Value = 1283.3334
$ws.Range("K89").Value = 6416.666999999999
$ws.Range("M89").Value = -800.6669999999995

# Row 105
$ws.Range("H105").Value = 1552.6428
$ws.Range("I105").Value = 1357.5
$ws.Range("K105").Value = 1357.5
$ws.Range("M105").Value = 389.5

# Row 133
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120

# Row 134
$ws.Range("H134").Value = 2808.2173
$ws.Range("I134").Value = 1717
$ws.Range("J134").Value = 5900
$ws.Range("K134").Value = 5151
$ws.Range("L134").Value = 17700
$ws.Range("M134").Value = -2616
$ws.Range("N134").Value = -22770

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2003061.9
$ws.Range("I31").Value = 2501987.2
$ws.Range("J31").Value = 7360
$ws.Range("K31").Value = 2501987.2
$ws.Range("L31").Value = 7360
$ws.Range("M31").Value = -2501692.2
$ws.Range("N31").Value = -7950

# Row 34
$ws.Range("H34").Value = 2003061.9
$ws.Range("I34").Value = 2501987.2
$ws.Range("J34").Value = 7360
$ws.Range("K34").Value = 2501987.2
$ws.Range("L34").Value = 7360
$ws.Range("M34").Value = -2501785.2
$ws.Range("N34").Value = -7764

# Row 127
$ws.Range("H127").Value = 39900
$ws.Range("J127").Value = 39900
$ws.Range("L127").Value = 39900
$ws.Range("N127").Value = -49820

# Row 134
$ws.Range("H134").Value = 2975.5
$ws.Range("I134").Value = 1413.7778
$ws.Range("K134").Value = 4241.3334
$ws.Range("M134").Value = -1706.3334

# Row 135
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

# Row 137
$ws.Range("H137").Value = 23582.223
$ws.Range("J137").Value = 23582.223
$ws.Range("L137").Value = 23582.223
$ws.Range("N137").Value = -33782.223

$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 2428.5715
$ws.Range("J58").Value = 2750
$ws.Range("L58").Value = 8250
$ws.Range("N58").Value = -8506

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 961.5
$ws.Range("I31").Value = 956
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 956
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = -664
$ws.Range("N31").Value = -1584

# Row 37
$ws.Range("H37").Value = 961.5
$ws.Range("I37").Value = 956
$ws.Range("J37").Value = 1000
$ws.Range("K37").Value = 956
$ws.Range("L37").Value = 1000
$ws.Range("M37").Value = -679
$ws.Range("N37").Value = -1554

# Row 80
$ws.Range("H80").Value = 3201.2
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 2668.6667
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 2668.6667
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -4664.6667

# Row 83
$ws.Range("H83").Value = 3201.2
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 2668.6667
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 13343.3335
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -23327.3335

# Row 97
$ws.Range("H97").Value = 1093.2593
$ws.Range("I97").Value = 785.9583
$ws.Range("J97").Value = 3551.6667
$ws.Range("K97").Value = 785.9583
$ws.Range("L97").Value = 3551.6667
$ws.Range("M97").Value = -289.9583
$ws.Range("N97").Value = -4543.6667

# Row 135
$ws.Range("H135").Value = 28998.334
$ws.Range("J135").Value = 28998.334
$ws.Range("L135").Value = 28998.334
$ws.Range("N135").Value = -39138.334

# Row 137
$ws.Range("H137").Value = 29573.076
$ws.Range("J137").Value = 29573.076
$ws.Range("L137").Value = 29573.076
$ws.Range("N137").Value = -39773.076

$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Range("H133").Value = 29487.5
$ws.Range("J133").Value = 29487.5
$ws.Range("L133").Value = 29487.5
$ws.Range("N133").Value = -34547.5

# Row 135
$ws.Range("H135").Value = 29043
$ws.Range("J135").Value = 29043
$ws.Range("L135").Value = 29043
$ws.Range("N135").Value = -39183

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1114.0416
$ws.Range("I113").Value = 539.625
$ws.Range("J113").Value = 2262.875
$ws.Range("K113").Value = 1618.875
$ws.Range("L113").Value = 6788.625
$ws.Range("M113").Value = 551.125
$ws.Range("N113").Value = -11128.625
